$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns B:E to C:F (right to left so data isn't clobbered
# before it's copied). Column F is brand new, so copy column E's formatting
# onto it before overwriting its value.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

for ($row = 1; $row -le 5; $row++) {
    $ws.Cells.Item($row, 6).Value = $ws.Cells.Item($row, 5).Value2
    $ws.Cells.Item($row, 5).Value = $ws.Cells.Item($row, 4).Value2
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 3).Value = $ws.Cells.Item($row, 2).Value2
}

# New column B: the "Unnamed: 0" index column pandas writes for an
# unnamed/duplicated index, mirroring column A's values.
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
